# Apply the competition-sheet edits described by the commit:
#   - row 2: record scores 2 / 3 as text "2.0" / "3.0"
#   - row 4: team renamed "大學男"->"大學南", opponents "男二"->"南二", "男四"->"男三"
#   - row 5: team renamed "大學男"->"大學南", opponents "男二"->"南二", "男三"->"南寺"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 / F2 need to be stored as literal text "2.0"/"3.0" (not numbers).
# Force text formatting first so Excel doesn't auto-coerce the value to a
# number, then drop back to the default "Normal" style so no extra
# formatting is left behind on the cell.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.0"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3.0"
$ws.Range("F2").Style = "Normal"

# Rows 4 and 5: team/opponent relabelling.
$ws.Range("B4").Value = "大學南"
$ws.Range("C4").Value = "南二"
$ws.Range("D4").Value = "男三"

$ws.Range("B5").Value = "大學南"
$ws.Range("C5").Value = "南二"
$ws.Range("D5").Value = "南寺"
